# Add team record columns (Wins/Losses/Ties) to the DET_2004 sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): AD1=Wins, AE1=Losses, AF1=Ties ---
$ws.Range("AD1").Value() = "Wins"
$ws.Range("AE1").Value() = "Losses"
$ws.Range("AF1").Value() = "Ties"

# Match the formatting used by the other header cells (e.g. AC1):
# bold font, thin border, centered horizontally, top vertical alignment.
$header = $ws.Range("AD1:AF1")
$header.Font.Bold = $true
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4160
$header.Borders.LineStyle = 1
$header.Borders.Weight = 2

# --- Data rows (2-41): constant team record values ---
$lastRow = 41
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Range("AD" + $r).Value() = 72
    $ws.Range("AE" + $r).Value() = 90
    $ws.Range("AF" + $r).Value() = 0
}
